$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> updated values (D: Fecha serial, J: Volumen, K: Precio minimo,
# L: Precio maximo, M: Precio promedio ponderado, O: Origen, P: Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 },
    @{ Row = 3;  D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 },
    @{ Row = 4;  D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 },
    @{ Row = 5;  D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 },
    @{ Row = 6;  D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 },
    @{ Row = 7;  D = 45113; J = 8;  K = 17000; L = 17000; M = 17000; O = "Provincia de Limarí"; P = 680 },
    @{ Row = 8;  D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 },
    @{ Row = 9;  D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 },
    @{ Row = 10; D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 },
    @{ Row = 11; D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 },
    @{ Row = 12; D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 },
    @{ Row = 13; D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 },
    @{ Row = 14; D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
